# Update Sheets via scheduled runner: refresh Leve profit calculations
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 5 (Leve Item ID 5503)
$ws.Range("H5").Value = 1415.3
$ws.Range("I5").Value = 1681.625
$ws.Range("J5").Value = 350
$ws.Range("K5").Value = 1681.625
$ws.Range("L5").Value = 350
$ws.Range("M5").Value = -1566.625
$ws.Range("N5").Value = -580
# row 17 (Leve Item ID 38956)
$ws.Range("H17").Value = 22750.666
$ws.Range("J17").Value = 22750.666
$ws.Range("L17").Value = 68251.99800000001
$ws.Range("N17").Value = -68587.99800000001
# row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 2281.9092
$ws.Range("I40").Value = 2462.625
$ws.Range("J40").Value = 1800
$ws.Range("K40").Value = 2462.625
$ws.Range("L40").Value = 1800
$ws.Range("M40").Value = -2287.625
$ws.Range("N40").Value = -2150
# row 52 (Leve Item ID 4567)
$ws.Range("H52").Value = 2457.1428
$ws.Range("I52").Value = 400
$ws.Range("J52").Value = 4000
$ws.Range("K52").Value = 1200
$ws.Range("L52").Value = 12000
$ws.Range("M52").Value = -1040
$ws.Range("N52").Value = -12320
# row 70 (Leve Item ID 12604)
$ws.Range("H70").Value = 44761
$ws.Range("I70").Value = 167716.67
$ws.Range("J70").Value = 1364.8823
$ws.Range("K70").Value = 503150.01
$ws.Range("L70").Value = 4094.6469
$ws.Range("M70").Value = -502880.01
$ws.Range("N70").Value = -4634.6469
# row 73 (Leve Item ID 12604)
$ws.Range("H73").Value = 44761
$ws.Range("I73").Value = 167716.67
$ws.Range("J73").Value = 1364.8823
$ws.Range("K73").Value = 503150.01
$ws.Range("L73").Value = 4094.6469
$ws.Range("M73").Value = -502214.01
$ws.Range("N73").Value = -5966.6469
# row 100 (Leve Item ID 19906)
$ws.Range("H100").Value = 3668.3125
$ws.Range("I100").Value = 3669.5715
$ws.Range("J100").Value = 3667.3333
$ws.Range("K100").Value = 3669.5715
$ws.Range("L100").Value = 3667.3333
$ws.Range("M100").Value = -3128.5715
$ws.Range("N100").Value = -4749.3333
# row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 6668995.5
$ws.Range("I116").Value = 15386511
$ws.Range("J116").Value = 2659.9412
$ws.Range("K116").Value = 15386511
$ws.Range("L116").Value = 2659.9412
$ws.Range("M116").Value = -15383069
$ws.Range("N116").Value = -9543.941200000001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 34 (Leve Item ID 2753)
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
# row 37 (Leve Item ID 3096)
$ws.Range("H37").Value = 26000
$ws.Range("J37").Value = 28000
$ws.Range("L37").Value = 28000
$ws.Range("N37").Value = -28546
# row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 1122.5405
$ws.Range("I74").Value = 1316.1666
$ws.Range("J74").Value = 765.0769
$ws.Range("K74").Value = 1316.1666
$ws.Range("L74").Value = 765.0769
$ws.Range("M74").Value = -442.1666
$ws.Range("N74").Value = -2513.0769
# row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 1122.5405
$ws.Range("I77").Value = 1316.1666
$ws.Range("J77").Value = 765.0769
$ws.Range("K77").Value = 6580.833000000001
$ws.Range("L77").Value = 3825.3845
$ws.Range("M77").Value = -2212.833000000001
$ws.Range("N77").Value = -12561.3845

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 2113.6924
$ws.Range("I99").Value = 1466.6666
$ws.Range("K99").Value = 1466.6666
$ws.Range("M99").Value = 31.33339999999998
# row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 3579.7917
$ws.Range("I134").Value = 3609.0667
$ws.Range("J134").Value = 3531
$ws.Range("K134").Value = 10827.2001
$ws.Range("L134").Value = 10593
$ws.Range("M134").Value = -8292.2001
$ws.Range("N134").Value = -15663

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 1450
$ws.Range("I31").Value = 1158.569
$ws.Range("J31").Value = 2576.8667
$ws.Range("K31").Value = 1158.569
$ws.Range("L31").Value = 2576.8667
$ws.Range("M31").Value = -863.569
$ws.Range("N31").Value = -3166.8667
# row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 1450
$ws.Range("I34").Value = 1158.569
$ws.Range("J34").Value = 2576.8667
$ws.Range("K34").Value = 1158.569
$ws.Range("L34").Value = 2576.8667
$ws.Range("M34").Value = -956.569
$ws.Range("N34").Value = -2980.8667

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 76 (Leve Item ID 12869)
$ws.Range("H76").Value = 3430
$ws.Range("I76").Value = 866.6667
$ws.Range("J76").Value = 4711.6665
$ws.Range("K76").Value = 2600.0001
$ws.Range("L76").Value = 14134.9995
$ws.Range("M76").Value = -2217.0001
$ws.Range("N76").Value = -14900.9995
# row 79 (Leve Item ID 12869)
$ws.Range("H79").Value = 3430
$ws.Range("I79").Value = 866.6667
$ws.Range("J79").Value = 4711.6665
$ws.Range("K79").Value = 2600.0001
$ws.Range("L79").Value = 14134.9995
$ws.Range("M79").Value = -1274.0001
$ws.Range("N79").Value = -16786.9995
# row 113 (Leve Item ID 27843)
$ws.Range("H113").Value = 597
$ws.Range("I113").Value = 527.94446
$ws.Range("J113").Value = 1011.3333
$ws.Range("K113").Value = 1583.83338
$ws.Range("L113").Value = 3033.9999
$ws.Range("M113").Value = 586.16662
$ws.Range("N113").Value = -7373.9999
# row 133 (Leve Item ID 44073)
$ws.Range("H133").Value = 5195.722
$ws.Range("I133").Value = 2566.2
$ws.Range("J133").Value = 6207.077
$ws.Range("K133").Value = 7698.599999999999
$ws.Range("L133").Value = 18621.231
$ws.Range("M133").Value = -2638.599999999999
$ws.Range("N133").Value = -28741.231
# row 134 (Leve Item ID 44074)
$ws.Range("H134").Value = 4454.1
$ws.Range("I134").Value = 2520.9092
$ws.Range("J134").Value = 5573.316
$ws.Range("K134").Value = 7562.7276
$ws.Range("L134").Value = 16719.948
$ws.Range("M134").Value = -2492.7276
$ws.Range("N134").Value = -26859.948

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 933.3333
$ws.Range("I22").Value = 900
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 900
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -605
$ws.Range("N22").Value = -1590
# row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 933.3333
$ws.Range("I27").Value = 900
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 900
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = -793
$ws.Range("N27").Value = -1214
# row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 1200
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 1600
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 1600
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -1976
# row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 738.75
$ws.Range("I55").Value = 449
$ws.Range("J55").Value = 835.3333
$ws.Range("K55").Value = 449
$ws.Range("L55").Value = 835.3333
$ws.Range("M55").Value = -276
$ws.Range("N55").Value = -1181.3333

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 40 (Leve Item ID 3601)
$ws.Range("H40").Value = 20000
$ws.Range("J40").Value = 20000
$ws.Range("L40").Value = 20000
$ws.Range("N40").Value = -20298
# row 62 (Leve Item ID 12589)
$ws.Range("H62").Value = 4433.3335
$ws.Range("I62").Value = 4062.5
$ws.Range("J62").Value = 4857.143
$ws.Range("K62").Value = 4062.5
$ws.Range("L62").Value = 4857.143
$ws.Range("M62").Value = -3438.5
$ws.Range("N62").Value = -6105.143
# row 65 (Leve Item ID 12589)
$ws.Range("H65").Value = 4433.3335
$ws.Range("I65").Value = 4062.5
$ws.Range("J65").Value = 4857.143
$ws.Range("K65").Value = 20312.5
$ws.Range("L65").Value = 24285.715
$ws.Range("M65").Value = -17192.5
$ws.Range("N65").Value = -30525.715
# row 96 (Leve Item ID 19977)
$ws.Range("H96").Value = 3146.5
$ws.Range("I96").Value = 2844.2
$ws.Range("J96").Value = 3448.8
$ws.Range("K96").Value = 2844.2
$ws.Range("L96").Value = 3448.8
$ws.Range("M96").Value = -1471.2
$ws.Range("N96").Value = -6194.8
